# Auto-generated script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "29.090.56"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.843.27"
$ws.Range("E3").Value = "  -0.90%  "
Set-TextValue "D4" "0.9984"
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "245.51"
$ws.Range("E5").Value = "  +1.57%  "
Set-TextValue "D6" "0.6963"
$ws.Range("E6").Value = "  -0.66%  "
Set-TextValue "D7" "0.9990"
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue "D8" "0.07706"
$ws.Range("E8").Value = "  -1.42%  "
Set-TextValue "D9" "0.3054"
$ws.Range("E9").Value = "  -1.73%  "
Set-TextValue "D10" "23.53"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  +0.33%  "
Set-TextValue "D12" "92.94"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.121"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.833.81"
$ws.Range("E14").Value = "  -1.43%  "
Set-TextValue "D15" "0.6840"
$ws.Range("E15").Value = "  -0.97%  "
Set-TextValue "D16" "6.634"
$ws.Range("E16").Value = "  +1.14%  "
Set-TextValue "D17" "0.000008284"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "29.082.45"
$ws.Range("E18").Value = "  -0.46%  "
Set-TextValue "D19" "242.04"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "2.078.07"
$ws.Range("E20").Value = "  -1.58%  "
Set-TextValue "D21" "12.76"
$ws.Range("E21").Value = "  -1.28%  "
Set-TextValue "D22" "0.9989"
$ws.Range("E22").Value = "  -0.12%  "
Set-TextValue "D23" "7.479"
$ws.Range("E23").Value = "  -1.55%  "
Set-TextValue "D24" "0.9994"
$ws.Range("E24").Value = "  -0.05%  "
Set-TextValue "D25" "0.1508"
$ws.Range("E25").Value = "  -1.63%  "
Set-TextValue "D26" "159.16"
$ws.Range("E26").Value = "  -1.06%  "
Set-TextValue "D27" "8.811"
$ws.Range("E27").Value = "  -0.89%  "
Set-TextValue "D28" "18.22"
$ws.Range("E28").Value = "  -1.77%  "
Set-TextValue "D29" "1.541"
$ws.Range("E29").Value = "  -1.96%  "
Set-TextValue "D30" "4.230"
$ws.Range("E30").Value = "  -1.11%  "
Set-TextValue "D31" "4.174"
$ws.Range("E31").Value = "  -1.78%  "
Set-TextValue "D32" "1.200"
$ws.Range("E32").Value = "  -1.16%  "
Set-TextValue "D33" "0.05108"
$ws.Range("E33").Value = "  -2.20%  "
Set-TextValue "D34" "0.7877"
$ws.Range("E34").Value = "  +4.02%  "
Set-TextValue "D35" "1.864"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  -2.57%  "
Set-TextValue "D37" "2.695"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "1.309.26"
$ws.Range("E38").Value = "  +7.20%  "
Set-TextValue "D39" "0.01863"
$ws.Range("E39").Value = "  +0.02%  "
Set-TextValue "D40" "2.705"
$ws.Range("E40").Value = "  -0.60%  "
Set-TextValue "D41" "0.9483"
$ws.Range("E41").Value = "  +5.15%  "
Set-TextValue "D42" "6.134"
$ws.Range("E42").Value = "  +5.39%  "
Set-TextValue "D43" "107.65"
$ws.Range("E43").Value = "  -2.26%  "
Set-TextValue "D44" "0.9991"
$ws.Range("E44").Value = "  -0.02%  "
Set-TextValue "D45" "9.721"
$ws.Range("E45").Value = "  +2.12%  "
Set-TextValue "D46" "0.5170"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "1.979.46"
$ws.Range("E47").Value = "  -1.51%  "
Set-TextValue "D48" "64.16"
$ws.Range("E48").Value = "  -2.42%  "
Set-TextValue "D49" "1.760"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -4.03%  "
Set-TextValue "D51" "6.985"
$ws.Range("E51").Value = "  -0.72%  "
